$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "run_count" column header (L1) to "processed"
$ws.Range("L1").Value = "processed"

# Update the active selection to L1 (was I9)
$ws.Range("L1").Select()
